$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.879.48"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "3.759.55"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "621.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "3.756.12"
$ws.Range("E7").Value = "  +2.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E10").Value = "  +3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.39%  "
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").Value = "4.379.35"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").Value = "3.757.80"
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").Value = "70.008.35"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "508.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.22%  "
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.23%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E28").Value = "  +21.04%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.337"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "426.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").Value = "3.009.83"
$ws.Range("E46").Value = "  -4.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0365"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.20%  "

Write-Output "Updated cryptos list"
